$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.582.92'
$ws.Range('E2').Value = '  +4.16%  '
$ws.Range('D3').Value = '3.486.49'
$ws.Range('E3').Value = '  +2.65%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'591.98"
$ws.Range('E5').Value = '  +3.88%  '
$ws.Range('D6').Value = "'169.13"
$ws.Range('E6').Value = '  +4.55%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.485.21'
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('D9').Value = "'0.594"
$ws.Range('E9').Value = '  +8.74%  '
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('E11').Value = '  +7.22%  '
$ws.Range('E12').Value = '  +4.45%  '
$ws.Range('D13').Value = '4.041.86'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').Value = "'28.12"
$ws.Range('E15').Value = '  +4.94%  '
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('D17').Value = '66.629.40'
$ws.Range('E17').Value = '  +4.20%  '
$ws.Range('D18').Value = '3.495.47'
$ws.Range('E18').Value = '  +2.75%  '
$ws.Range('E19').Value = '  +3.48%  '
$ws.Range('D20').Value = "'14.02"
$ws.Range('E20').Value = '  +4.24%  '
$ws.Range('D21').Value = "'392.25"
$ws.Range('E21').Value = '  +5.72%  '
$ws.Range('E22').Value = '  +2.05%  '
$ws.Range('D23').Value = "'73.01"
$ws.Range('E23').Value = '  +4.29%  '
$ws.Range('D24').Value = "'1.00"
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = "'0.535"
$ws.Range('E25').Value = '  +4.97%  '
$ws.Range('E26').Value = '  +6.40%  '
$ws.Range('D27').Value = "'10.33"
$ws.Range('E27').Value = '  +9.34%  '
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = "'6.32"
$ws.Range('E30').Value = '  +4.98%  '
$ws.Range('E31').Value = '  +5.85%  '
$ws.Range('D32').Value = "'2.07"
$ws.Range('E32').Value = '  +3.92%  '
$ws.Range('D33').Value = "'23.63"
$ws.Range('E33').Value = '  +4.26%  '
$ws.Range('D34').Value = "'7.40"
$ws.Range('E34').Value = '  +6.57%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = "'1.61"
$ws.Range('E36').Value = '  +10.03%  '
$ws.Range('D37').Value = "'161.36"
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('D38').Value = "'0.900"
$ws.Range('E38').Value = '  +4.77%  '
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('D40').Value = "'6.77"
$ws.Range('E40').Value = '  +5.95%  '
$ws.Range('D41').Value = "'0.0746"
$ws.Range('E41').Value = '  +4.03%  '
$ws.Range('D42').Value = "'26.56"
$ws.Range('E42').Value = '  +3.61%  '
$ws.Range('D43').Value = "'4.63"
$ws.Range('E43').Value = '  +7.07%  '
$ws.Range('D44').Value = "'26.73"
$ws.Range('E44').Value = '  +3.92%  '
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('D46').Value = '2.767.17'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').Value = "'0.0314"
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('D48').Value = "'2.49"
$ws.Range('E48').Value = '  +4.54%  '
$ws.Range('D49').Value = "'346.59"
$ws.Range('E49').Value = '  +5.41%  '
$ws.Range('E50').Value = '  +5.60%  '
$ws.Range('D51').Value = "'33.93"
$ws.Range('E51').Value = '  +13.19%  '
